$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 98
$ws.Range("H4").Value = 1142
$ws.Range("I4").Value = 1021
$ws.Range("J4").Value = 1106
$ws.Range("Q4").Value = 740
